$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; existing rows 57-88 shift down to 58-89.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new data point.
$ws.Range("A57").Value = 1
$ws.Range("B57").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C57").Value = "Arica y Parinacota"
$ws.Range("D57").Value = (Get-Date -Year 2022 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E57").Value = 15
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100106
$ws.Range("H57").Value = "Oleaginosos"
$ws.Range("I57").Value = 100106002
$ws.Range("J57").Value = "Palta"
$ws.Range("K57").Value = "Hass"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 400
$ws.Range("N57").Value = 16000
$ws.Range("O57").Value = 17000
$ws.Range("P57").Value = 16500
$ws.Range("Q57").Value = "$/bandeja 10 kilos"
$ws.Range("R57").Value = "Perú"
$ws.Range("S57").Value = 1650
$ws.Range("T57").Value = 10
